$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2 = @(0.087324523441786, 0.532791898334032, 0.5528062467454449, 0.743509412681134, 0.7610876864711806, 17)
    3 = @(0.001144284503009711, 0.3931116398949452, 0.2245968071122517, 0.4739164558361017, 0.4885004773109793, 17)
    4 = @(0.2183027114686655, 0.4424788170065723, 0.3666827254751553, 0.6055433307990067, 0.5822077098150801, 17)
    5 = @(0.07944540621654918, 0.4860981583435489, 0.3614563784931307, 0.6012124237681143, 0.6142811703890018, 17)
    6 = @(0.2020838518978592, 0.3911042730239319, 0.2849992644578546, 0.5338532237027839, 0.5093343377279087, 17)
    7 = @(0.148083212490263, 0.3930966733895834, 0.2550052339340027, 0.5049804292584047, 0.4976383057699656, 17)
    8 = @(0.2335018036174021, 0.4725336998297775, 0.4322111561822374, 0.6574276813325077, 0.6334773617760631, 17)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
